# Daily IST report: add CSV/MD/XLSX
# Adds the 2026-02-23 submission column, inserting it before the
# total_files / unique_days summary columns, and recomputes those
# summary columns for every contributor row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (pushes total_files -> J, unique_days -> K).
# Excel's native Insert copies formatting from the column to the left,
# which gives the new column the same header style as the other date
# columns (D1:H1) for free.
$ws.Columns("I:I").Insert()

# The inserted column defaults to the sheet's base width; the date
# columns use an explicit width of 12 (ColumnWidth 12 - 0.83 offset
# between the COM ColumnWidth property and the stored OOXML width).
$ws.Columns("I:I").ColumnWidth = 11.17

# New date header.
$ws.Range("I1").Value = "2026-02-23"

# Per-contributor submission flag for 2026-02-23, and the recomputed
# total_files / unique_days summary columns.
$newDayFlag = @(1,1,1,1,1,1,1,1,0,1,1,0,1,1,1,1,1,0,1,0,0,0,1,1,0,1,0,0,0,1,1,1,1,1,1,0,0,0,1,1,1,0,0,1,1,1,0,1,1,0,0,1,0,1,0,1,1,1,1,0,0,0,1,0,0,1,0,0,1,1,0,0,1,0,1,0,0,1,1,1,1,0,0,0,1,0,1,0,1,0,0,1,1,1,0,0,0,0,0,1,0,1,0,1,1,0,1,0)
$totalFiles = @(3,2,2,3,1,3,3,3,2,2,3,2,3,2,1,3,3,2,2,0,0,0,2,3,0,4,0,0,0,3,3,3,3,3,2,1,0,0,3,3,3,21,0,2,3,3,0,3,2,0,2,2,0,3,0,2,2,2,3,1,0,0,3,0,0,1,0,0,3,1,0,0,15,0,2,0,0,1,2,2,3,0,1,0,2,0,2,0,1,0,0,1,2,32,0,0,0,0,1,3,0,3,12,3,3,0,2,0)
$uniqueDays = @(3,2,2,3,1,3,3,3,2,2,3,2,3,2,1,3,3,2,2,0,0,0,2,3,0,2,0,0,0,3,3,3,3,3,2,1,0,0,3,3,3,1,0,2,3,3,0,3,2,0,2,2,0,3,0,2,2,2,3,1,0,0,3,0,0,1,0,0,3,1,0,0,3,0,2,0,0,1,2,2,3,0,1,0,2,0,2,0,1,0,0,1,2,3,0,0,0,0,1,3,0,3,1,3,3,0,2,0)

$rowCount = $newDayFlag.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $newDayFlag[$i]
    $ws.Cells.Item($r, 10).Value = $totalFiles[$i]
    $ws.Cells.Item($r, 11).Value = $uniqueDays[$i]
}

Write-Output "Applied 2026-02-23 column; used range is $($ws.UsedRange.Address())"
